$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.80835738505953203
$ws.Range("I3").Value = 0.59145229122494702
$ws.Range("I4").Value = 0.52639073649588597
$ws.Range("I7").Value = 2.5210337555198099
$ws.Range("I10").Value = 0.41737740556553699
$ws.Range("I13").Value = 0.54058373742891996
$ws.Range("I14").Value = 0.56066093329653799
$ws.Range("I16").Value = 0.68684069528846103
$ws.Range("I17").Value = 0.675652400825275

$ws.Range("J3").Select()
